# Auto-generated edit script for data_sigla.xlsx update v2
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BASELINE_SIGLA")
$ws2 = $wb.Worksheets.Item("Controle")

# ---- Sheet1: BASELINE_SIGLA (rows 2-15) ----
# Row 2: MA0290
$ws1.Range("A2").Value = "MA0290"
$ws1.Range("B2").Value = "Development"
$ws1.Range("C2").Value = "BACKUP 8GB RAM 8vCPU (BACKUP)"
$ws1.Range("D2").Value = 230
$ws1.Range("E2").Value = 8
$ws1.Range("F2").Value = 8
$ws1.Range("G2").Value = "BACKUP"
$ws1.Range("H2").Value = "Windows Server"

# Row 3: MA0291
$ws1.Range("A3").Value = "MA0291"
$ws1.Range("B3").Value = "Development"
$ws1.Range("C3").Value = "BACKUP 8GB RAM 8vCPU (BACKUP)"
$ws1.Range("D3").Value = 170
$ws1.Range("E3").Value = 8
$ws1.Range("F3").Value = 8
$ws1.Range("G3").Value = "BACKUP"
$ws1.Range("H3").Value = "Windows Server"

# Row 4: MA0292
$ws1.Range("A4").Value = "MA0292"
$ws1.Range("B4").Value = "Homologation"
$ws1.Range("C4").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D4").Value = 230
$ws1.Range("E4").Value = 16
$ws1.Range("F4").Value = 16
$ws1.Range("G4").Value = "WEBSERVER"
$ws1.Range("H4").Value = "Linux Server"

# Row 5: MA0293
$ws1.Range("A5").Value = "MA0293"
$ws1.Range("B5").Value = "Homologation"
$ws1.Range("C5").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D5").Value = 230
$ws1.Range("E5").Value = 16
$ws1.Range("F5").Value = 16
$ws1.Range("G5").Value = "WEBSERVER"
$ws1.Range("H5").Value = "Linux Server"

# Row 6: MA0294
$ws1.Range("A6").Value = "MA0294"
$ws1.Range("B6").Value = "Production"
$ws1.Range("C6").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D6").Value = 230
$ws1.Range("E6").Value = 16
$ws1.Range("F6").Value = 16
$ws1.Range("G6").Value = "WEBSERVER"
$ws1.Range("H6").Value = "Linux Server"

# Row 7: MA0295
$ws1.Range("A7").Value = "MA0295"
$ws1.Range("B7").Value = "Production"
$ws1.Range("C7").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D7").Value = 230
$ws1.Range("E7").Value = 16
$ws1.Range("F7").Value = 16
$ws1.Range("G7").Value = "WEBSERVER"
$ws1.Range("H7").Value = "Linux Server"

# Row 8: MA0296
$ws1.Range("A8").Value = "MA0296"
$ws1.Range("B8").Value = "Production"
$ws1.Range("C8").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D8").Value = 230
$ws1.Range("E8").Value = 16
$ws1.Range("F8").Value = 16
$ws1.Range("G8").Value = "WEBSERVER"
$ws1.Range("H8").Value = "Linux Server"

# Row 9: MA0297
$ws1.Range("A9").Value = "MA0297"
$ws1.Range("B9").Value = "Production"
$ws1.Range("C9").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D9").Value = 230
$ws1.Range("E9").Value = 16
$ws1.Range("F9").Value = 16
$ws1.Range("G9").Value = "WEBSERVER"
$ws1.Range("H9").Value = "Linux Server"

# Row 10: MA0298
$ws1.Range("A10").Value = "MA0298"
$ws1.Range("B10").Value = "Homologation"
$ws1.Range("C10").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D10").Value = 230
$ws1.Range("E10").Value = 16
$ws1.Range("F10").Value = 16
$ws1.Range("G10").Value = "WEBSERVER"
$ws1.Range("H10").Value = "Linux Server"

# Row 11: MA0299
$ws1.Range("A11").Value = "MA0299"
$ws1.Range("B11").Value = "Homologation"
$ws1.Range("C11").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D11").Value = 230
$ws1.Range("E11").Value = 16
$ws1.Range("F11").Value = 16
$ws1.Range("G11").Value = "WEBSERVER"
$ws1.Range("H11").Value = "Linux Server"

# Row 12: MA0300
$ws1.Range("A12").Value = "MA0300"
$ws1.Range("B12").Value = "Homologation"
$ws1.Range("C12").Value = "BACKUP 16GB RAM 16vCPU (WEBSERVER)"
$ws1.Range("D12").Value = 230
$ws1.Range("E12").Value = 16
$ws1.Range("F12").Value = 16
$ws1.Range("G12").Value = "WEBSERVER"
$ws1.Range("H12").Value = "Linux Server"

# Row 13: MA0301
$ws1.Range("A13").Value = "MA0301"
$ws1.Range("B13").Value = "Pre-Production"
$ws1.Range("C13").Value = "BACKUP 24GB RAM 16vCPU (BACKUP)"
$ws1.Range("D13").Value = 430
$ws1.Range("E13").Value = 24
$ws1.Range("F13").Value = 16
$ws1.Range("G13").Value = "BACKUP"
$ws1.Range("H13").Value = "Linux Server"

# Row 14: MA0302
$ws1.Range("A14").Value = "MA0302"
$ws1.Range("B14").Value = "Pre-Production"
$ws1.Range("C14").Value = "BACKUP 24GB RAM 16vCPU (BACKUP)"
$ws1.Range("D14").Value = 0
$ws1.Range("E14").Value = 24
$ws1.Range("F14").Value = 16
$ws1.Range("G14").Value = "BACKUP"
$ws1.Range("H14").Value = "Linux Server"

# Row 15: MA0303
$ws1.Range("A15").Value = "MA0303"
$ws1.Range("B15").Value = "Development"
$ws1.Range("C15").Value = "BACKUP 16GB RAM 16vCPU (BACKUP)"
$ws1.Range("D15").Value = 0
$ws1.Range("E15").Value = 16
$ws1.Range("F15").Value = 16
$ws1.Range("G15").Value = "BACKUP"
$ws1.Range("H15").Value = "Linux Server"

# ---- Sheet2: Controle (rows 2-6) ----
# Every column in this sheet stores its data as text (t="inlineStr"), even
# purely-numeric-looking values like instance counts/GB amounts. Plain
# `.Value = "10"` assignment would be auto-coerced to a number by Excel, so
# for those specific cells we first mark the cell as Text ("@") - reusing one
# shared style - to keep them typed the same way as the rest of the sheet.

# Row 2: Homologation / 10 x BACKUP 32GB RAM 16vCPU (WEBSERVER)
$ws2.Range("A2").Value = "Homologation"
$ws2.Range("B2").Value = "10 x BACKUP 32GB RAM 16vCPU (WEBSERVER)"
$ws2.Range("C2").Value = "sa-east-1"
$ws2.Range("D2").Value = "Linux"
$ws2.Range("E2").Value = "c6a.4xlarge"
$ws2.Range("F2").Value = "Shared Instances"
$ws2.Range("G2").NumberFormat = "@"
$ws2.Range("G2").Value = "10"
$ws2.Range("H2").NumberFormat = "@"
$ws2.Range("H2").Value = "40"
$ws2.Range("I2").Value = "Hours/Week"
$ws2.Range("J2").Value = "On-Demand"
$ws2.Range("K2").Value = "General Purpose SSD (gp3)"
$ws2.Range("L2").NumberFormat = "@"
$ws2.Range("L2").Value = "230"
$ws2.Range("O2").Value = "2x Daily"
$ws2.Range("P2").NumberFormat = "@"
$ws2.Range("P2").Value = "10"

# Row 3: Production / 12 x BACKUP 32GB RAM 16vCPU (WEBSERVER)
$ws2.Range("A3").Value = "Production"
$ws2.Range("B3").Value = "12 x BACKUP 32GB RAM 16vCPU (WEBSERVER)"
$ws2.Range("C3").Value = "sa-east-1"
$ws2.Range("D3").Value = "Linux"
$ws2.Range("E3").Value = "c6a.4xlarge"
$ws2.Range("F3").Value = "Shared Instances"
$ws2.Range("G3").NumberFormat = "@"
$ws2.Range("G3").Value = "12"
$ws2.Range("I3").Value = "Always On"
$ws2.Range("J3").Value = "1 Yr No Upfront EC2 Instance Savings Plan"
$ws2.Range("K3").Value = "General Purpose SSD (gp3)"
$ws2.Range("L3").NumberFormat = "@"
$ws2.Range("L3").Value = "230"
$ws2.Range("O3").Value = "6x Daily"
$ws2.Range("P3").NumberFormat = "@"
$ws2.Range("P3").Value = "20"

# Row 4: Development / 4 x BACKUP 16GB RAM 8vCPU (BACKUP)
$ws2.Range("A4").Value = "Development"
$ws2.Range("B4").Value = "4 x BACKUP 16GB RAM 8vCPU (BACKUP)"
$ws2.Range("C4").Value = "sa-east-1"
$ws2.Range("D4").Value = "Windows Server"
$ws2.Range("E4").Value = "c6i.2xlarge"
$ws2.Range("F4").Value = "Shared Instances"
$ws2.Range("G4").NumberFormat = "@"
$ws2.Range("G4").Value = "4"
$ws2.Range("H4").NumberFormat = "@"
$ws2.Range("H4").Value = "40"
$ws2.Range("I4").Value = "Hours/Week"
$ws2.Range("J4").Value = "On-Demand"
$ws2.Range("K4").Value = "General Purpose SSD (gp3)"
$ws2.Range("L4").NumberFormat = "@"
$ws2.Range("L4").Value = "230"
$ws2.Range("O4").Value = "2x Daily"
$ws2.Range("P4").NumberFormat = "@"
$ws2.Range("P4").Value = "10"

# Row 5: Pre-Production / 4 x BACKUP 32GB RAM 16vCPU (BACKUP)
$ws2.Range("A5").Value = "Pre-Production"
$ws2.Range("B5").Value = "4 x BACKUP 32GB RAM 16vCPU (BACKUP)"
$ws2.Range("C5").Value = "sa-east-1"
$ws2.Range("D5").Value = "Linux"
$ws2.Range("E5").Value = "c6a.4xlarge"
$ws2.Range("F5").Value = "Shared Instances"
$ws2.Range("G5").NumberFormat = "@"
$ws2.Range("G5").Value = "4"
$ws2.Range("H5").NumberFormat = "@"
$ws2.Range("H5").Value = "40"
$ws2.Range("I5").Value = "Hours/Week"
$ws2.Range("J5").Value = "On-Demand"
$ws2.Range("K5").Value = "General Purpose SSD (gp3)"
$ws2.Range("L5").NumberFormat = "@"
$ws2.Range("L5").Value = "430"
$ws2.Range("O5").Value = "2x Daily"
$ws2.Range("P5").NumberFormat = "@"
$ws2.Range("P5").Value = "10"

# Row 6: Development / 2 x BACKUP 32GB RAM 16vCPU (BACKUP)
$ws2.Range("A6").Value = "Development"
$ws2.Range("B6").Value = "2 x BACKUP 32GB RAM 16vCPU (BACKUP)"
$ws2.Range("C6").Value = "sa-east-1"
$ws2.Range("D6").Value = "Linux"
$ws2.Range("E6").Value = "c6a.4xlarge"
$ws2.Range("F6").Value = "Shared Instances"
$ws2.Range("G6").NumberFormat = "@"
$ws2.Range("G6").Value = "2"
$ws2.Range("H6").NumberFormat = "@"
$ws2.Range("H6").Value = "40"
$ws2.Range("I6").Value = "Hours/Week"
$ws2.Range("J6").Value = "On-Demand"
$ws2.Range("K6").Value = "General Purpose SSD (gp3)"
$ws2.Range("L6").NumberFormat = "@"
$ws2.Range("L6").Value = "0"
$ws2.Range("O6").Value = "2x Daily"
$ws2.Range("P6").NumberFormat = "@"
$ws2.Range("P6").Value = "10"

